$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells M1, N1 (match the bold/centered/bordered header style used by A1:L1)
$ws.Cells.Item(1, 13).Value = "respuesta"
$ws.Cells.Item(1, 14).Value = "razon"
$ws.Range("M1:N1").Font.Bold = $true
$ws.Range("M1:N1").HorizontalAlignment = -4108
$ws.Range("M1:N1").VerticalAlignment = -4160
$ws.Range("M1:N1").Borders.LineStyle = 1

# Update row 2 existing values
$ws.Cells.Item(2, 6).Value = " F, C'"
$ws.Cells.Item(2, 9).Value = " A"
$ws.Cells.Item(2, 13).Value = "un bicho"
$ws.Cells.Item(2, 14).Value = "es gris"

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "1"
$ws.Cells.Item(3, 2).Style = $ws.Cells.Item(2, 2).Style
$ws.Cells.Item(3, 3).Value = "?"
$ws.Cells.Item(3, 4).Value = "?"
$ws.Cells.Item(3, 5).Value = "?"
$ws.Cells.Item(3, 6).Value = " F, C"
$ws.Cells.Item(3, 7).Value = "?"
$ws.Cells.Item(3, 8).Value = "?"
$ws.Cells.Item(3, 9).Value = " Fi"
$ws.Cells.Item(3, 10).Value = "?"
$ws.Cells.Item(3, 11).Value = "?"
$ws.Cells.Item(3, 12).Value = "?"
$ws.Cells.Item(3, 13).Value = "un fuego"
$ws.Cells.Item(3, 14).Value = "es rojo"

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "1"
$ws.Cells.Item(4, 2).Style = $ws.Cells.Item(2, 2).Style
$ws.Cells.Item(4, 3).Value = "?"
$ws.Cells.Item(4, 4).Value = "?"
$ws.Cells.Item(4, 5).Value = "?"
$ws.Cells.Item(4, 6).Value = " F"
$ws.Cells.Item(4, 7).Value = "?"
$ws.Cells.Item(4, 8).Value = "?"
$ws.Cells.Item(4, 9).Value = " H, Hd"
$ws.Cells.Item(4, 10).Value = "Po3"
$ws.Cells.Item(4, 11).Value = "?"
$ws.Cells.Item(4, 12).Value = "?"
$ws.Cells.Item(4, 13).Value = "personas"
$ws.Cells.Item(4, 14).Value = "piernas"
